$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 119 (rows 119-176 shift down to 120-177)
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new data record
$ws.Cells.Item(119, 1).Value = 10
$ws.Cells.Item(119, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(119, 3).Value = "La Araucanía"
$ws.Cells.Item(119, 4).Value = 44960
$ws.Cells.Item(119, 5).Value = 9
$ws.Cells.Item(119, 6).Value = 100112031
$ws.Cells.Item(119, 7).Value = "Poroto verde"
$ws.Cells.Item(119, 8).Value = "Brío"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 50
$ws.Cells.Item(119, 11).Value = 1200
$ws.Cells.Item(119, 12).Value = 1200
$ws.Cells.Item(119, 13).Value = 1200
$ws.Cells.Item(119, 14).Value = "$/kilo"
$ws.Cells.Item(119, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(119, 16).Value = 1200
$ws.Cells.Item(119, 17).Value = 1
$ws.Cells.Item(119, 18).Value = "Hortaliza"
